# Mixer1: BOM - wip
# Fill in the "already have" (C) column for the parts rows that were
# entered on this pass (matching the B quantities so the E "still need"
# formula -> 0), then leave the current selection on C11 where editing
# continued.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 7
$ws.Range("C17").Value = 1
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 4

$ws.Range("C11").Select() | Out-Null
